$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.445.22"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").Value = "3.490.10"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.67"
$ws.Range("E5").Value = "  +5.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.28"
$ws.Range("E6").Value = "  +1.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -3.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.649"
$ws.Range("E10").Value = "  +0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.80"
$ws.Range("E11").Value = "  -2.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000305"
$ws.Range("E12").Value = "  -3.94%  "

$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").Value = "4.049.04"
$ws.Range("E14").Value = "  -1.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "609.57"
$ws.Range("E15").Value = "  +6.27%  "

$ws.Range("D16").Value = "69.465.10"
$ws.Range("E16").Value = "  -1.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.86"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.50"
$ws.Range("E18").Value = "  -1.89%  "

$ws.Range("D19").Value = "3.491.92"
$ws.Range("E19").Value = "  -1.81%  "

$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("E21").Value = "  -1.82%  "

$ws.Range("E22").Value = "  -3.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.28"
$ws.Range("E23").Value = "  +11.61%  "

$ws.Range("E24").Value = "  +3.02%  "

$ws.Range("E25").Value = "  +4.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.02"
$ws.Range("E26").Value = "  +3.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  -1.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  +4.22%  "

$ws.Range("E29").Value = "  +2.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.92"
$ws.Range("E30").Value = "  -3.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.57"
$ws.Range("E31").Value = "  +2.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.96"
$ws.Range("E32").Value = "  +5.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.115"
$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.25"
$ws.Range("E34").Value = "  -0.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.12"
$ws.Range("E35").Value = "  -5.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.66"
$ws.Range("E37").Value = "  +6.11%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "3.627.93"
$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("E39").Value = "  -4.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "502.85"
$ws.Range("E40").Value = "  -5.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.51"
$ws.Range("E41").Value = "  -4.48%  "

$ws.Range("D42").Value = "0.0₃0772"
$ws.Range("E42").Value = "  -4.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  -2.25%  "

$ws.Range("E44").Value = "  -1.93%  "

$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("E46").Value = "  +2.52%  "

$ws.Range("E47").Value = "  -4.48%  "

$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("E49").Value = "  -6.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.98"
$ws.Range("E50").Value = "  -3.46%  "

$ws.Range("E51").Value = "  -7.36%  "
